$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix quantities (digitation error correction)
$ws.Range("C5").Value = 3
$ws.Range("C24").Value = 3

# Row 25 used to be "Triple Tactic Talent" / "Non engine" -> replace with
# the card that used to live two rows below ("One for one" / "Starter")
$ws.Range("A25").Value = "One for one"
$ws.Range("B25").Value = "Starter"

# Row 26 used to be "One for one" / "Starter" -> replace with the card that
# used to live at the very bottom of the list ("Divine Temple of the Snale-Eye" / "Garnet")
$ws.Range("A26").Value = "Divine Temple of the Snale-Eye"
$ws.Range("B26").Value = "Garnet"

# The two trailing duplicate rows are no longer needed, so delete them
# (this shifts nothing further up since they are the last rows)
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()

# Restore the last active-cell selection recorded in the sheet view
$ws.Range("L7").Select()
